$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.787.19"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.220.41"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'241.66"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'72.84"
$ws.Range("E7").Value = "  -5.68%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.597"
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("D10").Value = "'42.25"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'6.96"
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "2.551.02"
$ws.Range("D15").Value = "'14.28"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "2.218.06"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "41.695.37"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'0.0000105"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("D20").Value = "'72.65"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'6.17"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'11.05"
$ws.Range("E22").Value = "  +20.07%  "
$ws.Range("D23").Value = "'229.53"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("E24").Value = "  -8.17%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'11.37"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "'166.97"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'20.48"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'0.0797"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").Value = "'5.54"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("D34").Value = "'30.18"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  -10.15%  "
$ws.Range("D37").Value = "'4.29"
$ws.Range("E37").Value = "  -6.35%  "
$ws.Range("D38").Value = "'0.0303"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").Value = "'13.36"
$ws.Range("E39").Value = "  -7.85%  "
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'64.67"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "'5.62"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").Value = "'8.72"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'103.29"
$ws.Range("E45").Value = "  -4.71%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "2.423.87"
$ws.Range("E51").Value = "  -1.72%  "
